# Apply updated TPM-derived values (NATMI Efnb1-Ephb4 LR-pair output) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = [double]"10.46510533333333"
$ws.Cells.Item(2, 8).Value = [double]"31.395316"
$ws.Cells.Item(2, 9).Value = [double]"0.5554075997074865"
$ws.Cells.Item(2, 10).Value = [double]"0.5554075997074865"
$ws.Cells.Item(2, 13).Value = [double]"34.82273866666667"
$ws.Cells.Item(2, 14).Value = [double]"104.468216"
$ws.Cells.Item(2, 15).Value = [double]"0.7026863693863559"
$ws.Cells.Item(2, 16).Value = [double]"0.702686369386356"
$ws.Cells.Item(2, 17).Value = [double]"364.4236281418062"
$ws.Cells.Item(2, 18).Value = [double]"3279.812653276256"
$ws.Cells.Item(2, 19).Value = [double]"0.3902773497680442"
$ws.Cells.Item(2, 20).Value = [double]"0.3902773497680442"

# Row 3
$ws.Cells.Item(3, 7).Value = [double]"10.46510533333333"
$ws.Cells.Item(3, 8).Value = [double]"31.395316"
$ws.Cells.Item(3, 9).Value = [double]"0.5554075997074865"
$ws.Cells.Item(3, 10).Value = [double]"0.5554075997074865"
$ws.Cells.Item(3, 15).Value = [double]"0.1722322446965897"
$ws.Cells.Item(3, 16).Value = [double]"0.1722322446965897"
$ws.Cells.Item(3, 17).Value = [double]"89.32221006385913"
$ws.Cells.Item(3, 18).Value = [double]"803.8998905747321"
$ws.Cells.Item(3, 19).Value = [double]"0.09565909761916537"
$ws.Cells.Item(3, 20).Value = [double]"0.09565909761916538"

# Row 4
$ws.Cells.Item(4, 7).Value = [double]"10.46510533333333"
$ws.Cells.Item(4, 8).Value = [double]"31.395316"
$ws.Cells.Item(4, 9).Value = [double]"0.5554075997074865"
$ws.Cells.Item(4, 10).Value = [double]"0.5554075997074865"
$ws.Cells.Item(4, 13).Value = [double]"6.169951999999999"
$ws.Cells.Item(4, 14).Value = [double]"18.509856"
$ws.Cells.Item(4, 15).Value = [double]"0.1245031647760143"
$ws.Cells.Item(4, 16).Value = [double]"0.1245031647760143"
$ws.Cells.Item(4, 17).Value = [double]"64.56919758161067"
$ws.Cells.Item(4, 18).Value = [double]"581.122778234496"
$ws.Cells.Item(4, 19).Value = [double]"0.06915000390423179"
$ws.Cells.Item(4, 20).Value = [double]"0.06915000390423182"

# Row 5
$ws.Cells.Item(5, 7).Value = [double]"10.46510533333333"
$ws.Cells.Item(5, 8).Value = [double]"31.395316"
$ws.Cells.Item(5, 9).Value = [double]"0.5554075997074865"
$ws.Cells.Item(5, 10).Value = [double]"0.5554075997074865"
$ws.Cells.Item(5, 13).Value = [double]"0.02865466666666667"
$ws.Cells.Item(5, 14).Value = [double]"0.085964"
$ws.Cells.Item(5, 15).Value = [double]"0.000578221141039957"
$ws.Cells.Item(5, 16).Value = [double]"0.0005782211410399571"
$ws.Cells.Item(5, 17).Value = [double]"0.2998741049582222"
$ws.Cells.Item(5, 18).Value = [double]"2.698866944624"
$ws.Cells.Item(5, 19).Value = [double]"0.0003211484160451266"
$ws.Cells.Item(5, 20).Value = [double]"0.0003211484160451266"

# Row 6
$ws.Cells.Item(6, 9).Value = [double]"0.3053945925621632"
$ws.Cells.Item(6, 10).Value = [double]"0.3053945925621632"
$ws.Cells.Item(6, 13).Value = [double]"34.82273866666667"
$ws.Cells.Item(6, 14).Value = [double]"104.468216"
$ws.Cells.Item(6, 15).Value = [double]"0.7026863693863559"
$ws.Cells.Item(6, 16).Value = [double]"0.702686369386356"
$ws.Cells.Item(6, 17).Value = [double]"200.3807752990889"
$ws.Cells.Item(6, 18).Value = [double]"1803.4269776918"
$ws.Cells.Item(6, 19).Value = [double]"0.2145966174777318"
$ws.Cells.Item(6, 20).Value = [double]"0.2145966174777319"

# Row 7
$ws.Cells.Item(7, 9).Value = [double]"0.3053945925621632"
$ws.Cells.Item(7, 10).Value = [double]"0.3053945925621632"
$ws.Cells.Item(7, 15).Value = [double]"0.1722322446965897"
$ws.Cells.Item(7, 16).Value = [double]"0.1722322446965897"
$ws.Cells.Item(7, 19).Value = [double]"0.0525987961951818"
$ws.Cells.Item(7, 20).Value = [double]"0.05259879619518181"

# Row 8
$ws.Cells.Item(8, 9).Value = [double]"0.3053945925621632"
$ws.Cells.Item(8, 10).Value = [double]"0.3053945925621632"
$ws.Cells.Item(8, 13).Value = [double]"6.169951999999999"
$ws.Cells.Item(8, 14).Value = [double]"18.509856"
$ws.Cells.Item(8, 15).Value = [double]"0.1245031647760143"
$ws.Cells.Item(8, 16).Value = [double]"0.1245031647760143"
$ws.Cells.Item(8, 17).Value = [double]"35.50380620986667"
$ws.Cells.Item(8, 18).Value = [double]"319.5342558888"
$ws.Cells.Item(8, 19).Value = [double]"0.03802259327947075"
$ws.Cells.Item(8, 20).Value = [double]"0.03802259327947077"

# Row 9
$ws.Cells.Item(9, 9).Value = [double]"0.3053945925621632"
$ws.Cells.Item(9, 10).Value = [double]"0.3053945925621632"
$ws.Cells.Item(9, 13).Value = [double]"0.02865466666666667"
$ws.Cells.Item(9, 14).Value = [double]"0.085964"
$ws.Cells.Item(9, 15).Value = [double]"0.000578221141039957"
$ws.Cells.Item(9, 16).Value = [double]"0.0005782211410399571"
$ws.Cells.Item(9, 17).Value = [double]"0.1648877871888889"
$ws.Cells.Item(9, 18).Value = [double]"1.4839900847"
$ws.Cells.Item(9, 19).Value = [double]"0.0001765856097787268"
$ws.Cells.Item(9, 20).Value = [double]"0.0001765856097787268"

# Row 10
$ws.Cells.Item(10, 7).Value = [double]"2.146766"
$ws.Cells.Item(10, 8).Value = [double]"6.440298"
$ws.Cells.Item(10, 9).Value = [double]"0.1139338891693565"
$ws.Cells.Item(10, 10).Value = [double]"0.1139338891693565"
$ws.Cells.Item(10, 13).Value = [double]"34.82273866666667"
$ws.Cells.Item(10, 14).Value = [double]"104.468216"
$ws.Cells.Item(10, 15).Value = [double]"0.7026863693863559"
$ws.Cells.Item(10, 16).Value = [double]"0.702686369386356"
$ws.Cells.Item(10, 17).Value = [double]"74.75627139648533"
$ws.Cells.Item(10, 18).Value = [double]"672.806442568368"
$ws.Cells.Item(10, 19).Value = [double]"0.08005979093048259"
$ws.Cells.Item(10, 20).Value = [double]"0.08005979093048261"

# Row 11
$ws.Cells.Item(11, 7).Value = [double]"2.146766"
$ws.Cells.Item(11, 8).Value = [double]"6.440298"
$ws.Cells.Item(11, 9).Value = [double]"0.1139338891693565"
$ws.Cells.Item(11, 10).Value = [double]"0.1139338891693565"
$ws.Cells.Item(11, 15).Value = [double]"0.1722322446965897"
$ws.Cells.Item(11, 16).Value = [double]"0.1722322446965897"
$ws.Cells.Item(11, 17).Value = [double]"18.32316804296067"
$ws.Cells.Item(11, 18).Value = [double]"164.908512386646"
$ws.Cells.Item(11, 19).Value = [double]"0.01962308947865074"
$ws.Cells.Item(11, 20).Value = [double]"0.01962308947865075"

# Row 12
$ws.Cells.Item(12, 7).Value = [double]"2.146766"
$ws.Cells.Item(12, 8).Value = [double]"6.440298"
$ws.Cells.Item(12, 9).Value = [double]"0.1139338891693565"
$ws.Cells.Item(12, 10).Value = [double]"0.1139338891693565"
$ws.Cells.Item(12, 13).Value = [double]"6.169951999999999"
$ws.Cells.Item(12, 14).Value = [double]"18.509856"
$ws.Cells.Item(12, 15).Value = [double]"0.1245031647760143"
$ws.Cells.Item(12, 16).Value = [double]"0.1245031647760143"
$ws.Cells.Item(12, 17).Value = [double]"13.245443175232"
$ws.Cells.Item(12, 18).Value = [double]"119.208988577088"
$ws.Cells.Item(12, 19).Value = [double]"0.01418512977682455"
$ws.Cells.Item(12, 20).Value = [double]"0.01418512977682455"

# Row 13
$ws.Cells.Item(13, 7).Value = [double]"2.146766"
$ws.Cells.Item(13, 8).Value = [double]"6.440298"
$ws.Cells.Item(13, 9).Value = [double]"0.1139338891693565"
$ws.Cells.Item(13, 10).Value = [double]"0.1139338891693565"
$ws.Cells.Item(13, 13).Value = [double]"0.02865466666666667"
$ws.Cells.Item(13, 14).Value = [double]"0.085964"
$ws.Cells.Item(13, 15).Value = [double]"0.000578221141039957"
$ws.Cells.Item(13, 16).Value = [double]"0.0005782211410399571"
$ws.Cells.Item(13, 17).Value = [double]"0.06151486414133333"
$ws.Cells.Item(13, 18).Value = [double]"0.553633777272"
$ws.Cells.Item(13, 19).Value = [double]"6.587898339862533E-05"
$ws.Cells.Item(13, 20).Value = [double]"6.587898339862536E-05"

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = [double]"0.4760280000000001"
$ws.Cells.Item(14, 8).Value = [double]"1.428084"
$ws.Cells.Item(14, 9).Value = [double]"0.02526391856099382"
$ws.Cells.Item(14, 10).Value = [double]"0.02526391856099382"
$ws.Cells.Item(14, 13).Value = [double]"34.82273866666667"
$ws.Cells.Item(14, 14).Value = [double]"104.468216"
$ws.Cells.Item(14, 15).Value = [double]"0.7026863693863559"
$ws.Cells.Item(14, 16).Value = [double]"0.702686369386356"
$ws.Cells.Item(14, 17).Value = [double]"16.576598642016"
$ws.Cells.Item(14, 18).Value = [double]"149.189387778144"
$ws.Cells.Item(14, 19).Value = [double]"0.01775261121009732"
$ws.Cells.Item(14, 20).Value = [double]"0.01775261121009732"

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = [double]"0.4760280000000001"
$ws.Cells.Item(15, 8).Value = [double]"1.428084"
$ws.Cells.Item(15, 9).Value = [double]"0.02526391856099382"
$ws.Cells.Item(15, 10).Value = [double]"0.02526391856099382"
$ws.Cells.Item(15, 15).Value = [double]"0.1722322446965897"
$ws.Cells.Item(15, 16).Value = [double]"0.1722322446965897"
$ws.Cells.Item(15, 17).Value = [double]"4.063014337452001"
$ws.Cells.Item(15, 18).Value = [double]"36.56712903706801"
$ws.Cells.Item(15, 19).Value = [double]"0.004351261403591803"
$ws.Cells.Item(15, 20).Value = [double]"0.004351261403591803"

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = [double]"0.4760280000000001"
$ws.Cells.Item(16, 8).Value = [double]"1.428084"
$ws.Cells.Item(16, 9).Value = [double]"0.02526391856099382"
$ws.Cells.Item(16, 10).Value = [double]"0.02526391856099382"
$ws.Cells.Item(16, 13).Value = [double]"6.169951999999999"
$ws.Cells.Item(16, 14).Value = [double]"18.509856"
$ws.Cells.Item(16, 15).Value = [double]"0.1245031647760143"
$ws.Cells.Item(16, 16).Value = [double]"0.1245031647760143"
$ws.Cells.Item(16, 17).Value = [double]"2.937069910656"
$ws.Cells.Item(16, 18).Value = [double]"26.433629195904"
$ws.Cells.Item(16, 19).Value = [double]"0.00314543781548722"
$ws.Cells.Item(16, 20).Value = [double]"0.003145437815487221"

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = [double]"0.4760280000000001"
$ws.Cells.Item(17, 8).Value = [double]"1.428084"
$ws.Cells.Item(17, 9).Value = [double]"0.02526391856099382"
$ws.Cells.Item(17, 10).Value = [double]"0.02526391856099382"
$ws.Cells.Item(17, 13).Value = [double]"0.02865466666666667"
$ws.Cells.Item(17, 14).Value = [double]"0.085964"
$ws.Cells.Item(17, 15).Value = [double]"0.000578221141039957"
$ws.Cells.Item(17, 16).Value = [double]"0.0005782211410399571"
$ws.Cells.Item(17, 17).Value = [double]"0.013640423664"
$ws.Cells.Item(17, 18).Value = [double]"0.122763812976"
$ws.Cells.Item(17, 19).Value = [double]"1.46081318174784E-05"
$ws.Cells.Item(17, 20).Value = [double]"1.46081318174784E-05"

